$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.00", "0.999") that Excel's
# COM layer would otherwise silently coerce to a Number on assignment
# (dropping formatting / trailing zeros). Force the Price column to Text
# first so every assignment below is stored as a literal inline/shared
# string, matching the workbook's original inlineStr cell type.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "72.122.91"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "4.038.46"
$ws.Range("E3").Value = "  +3.56%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "518.00"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "146.83"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").Value = "0.725"
$ws.Range("E7").Value = "  +18.62%  "
$ws.Range("D8").Value = "4.030.18"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +7.67%  "
$ws.Range("D11").Value = "0.175"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "0.0000327"
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("D13").Value = "47.39"
$ws.Range("E13").Value = "  +12.46%  "
$ws.Range("D14").Value = "11.07"
$ws.Range("E14").Value = "  +7.92%  "
$ws.Range("D15").Value = "4.699.60"
$ws.Range("E15").Value = "  +4.03%  "
$ws.Range("D16").Value = "4.062.45"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("D17").Value = "21.12"
$ws.Range("E17").Value = "  +6.53%  "
$ws.Range("D18").Value = "14.08"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "72.310.96"
$ws.Range("E21").Value = "  +4.44%  "
$ws.Range("D22").Value = "442.26"
$ws.Range("E22").Value = "  +3.74%  "
$ws.Range("D23").Value = "104.25"
$ws.Range("E23").Value = "  +18.24%  "
$ws.Range("D24").Value = "3.54"
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("D25").Value = "14.80"
$ws.Range("E25").Value = "  +4.60%  "
$ws.Range("D26").Value = "4.00"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "11.40"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "10.98"
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("D29").Value = "37.68"
$ws.Range("E29").Value = "  +3.54%  "
$ws.Range("D30").Value = "5.81"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("D31").Value = "3.24"
$ws.Range("E31").Value = "  +15.33%  "
$ws.Range("D32").Value = "13.66"
$ws.Range("E32").Value = "  +3.56%  "
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").Value = "682.77"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Value = "6.81"
$ws.Range("E35").Value = "  +14.52%  "
$ws.Range("D36").Value = "66.78"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("D37").Value = "42.52"
$ws.Range("E37").Value = "  +6.32%  "
$ws.Range("D38").Value = "0.0₃0858"
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").Value = "0.427"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").Value = "3.53"
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  +3.79%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "0.158"
$ws.Range("E46").Value = "  +12.40%  "
$ws.Range("D47").Value = "3.54"
$ws.Range("E47").Value = "  +3.69%  "
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("D49").Value = "3.05"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "9.14"
$ws.Range("E50").Value = "  +7.15%  "
$ws.Range("D51").Value = "3.32"
$ws.Range("E51").Value = "  +1.77%  "

# Restore the default (unstyled) cell style now that the values are safely
# stored as text, so the cells keep their original appearance/format.
$priceRange.Style = "Normal"
